$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 411; this shifts existing rows 411:504 down to 412:505
# and expands the used range to A1:R505 automatically.
$ws.Rows(411).Insert()

# Populate the newly inserted row 411 with the new record (a new weekly
# price observation). Fields A,B,C,E,F,G,H,I,N,O,Q,R mirror the record that
# used to sit at row 411 (now at row 412); D,J,K,L,M,P carry the new values.
$ws.Range("A411").Value = 8
$ws.Range("B411").Value = "Terminal La Palmera de La Serena"
$ws.Range("C411").Value = "Coquimbo"
$ws.Range("D411").Value = 45135
$ws.Range("E411").Value = 4
$ws.Range("F411").Value = 100112003
$ws.Range("G411").Value = "Ajo"
$ws.Range("H411").Value = "Chino"
$ws.Range("I411").Value = "Primera"
$ws.Range("J411").Value = 300
$ws.Range("K411").Value = 19500
$ws.Range("L411").Value = 20000
$ws.Range("M411").Value = 19750
$ws.Range("N411").Value = "$/caja 10 kilos"
$ws.Range("O411").Value = "China"
$ws.Range("P411").Value = 1975
$ws.Range("Q411").Value = 10
$ws.Range("R411").Value = "Hortaliza"
